$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 320, pushing the existing rows 320-337 down to 321-338.
$ws.Rows.Item(320).Insert()

# Populate the newly inserted row 320 with the new weekly price record.
$row = 320
$ws.Cells.Item($row,1).Value2  = 8
$ws.Cells.Item($row,2).Value2  = 'Terminal La Palmera de La Serena'
$ws.Cells.Item($row,3).Value2  = 'Coquimbo'
$ws.Cells.Item($row,4).Value2  = 45267
$ws.Cells.Item($row,4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item($row,5).Value2  = 4
$ws.Cells.Item($row,6).Value2  = 100112001
$ws.Cells.Item($row,7).Value2  = 'Berenjena'
$ws.Cells.Item($row,8).Value2  = 'Sin especificar'
$ws.Cells.Item($row,9).Value2  = 'Primera'
$ws.Cells.Item($row,10).Value2 = 500
$ws.Cells.Item($row,11).Value2 = 11000
$ws.Cells.Item($row,12).Value2 = 12000
$ws.Cells.Item($row,13).Value2 = 11500
$ws.Cells.Item($row,14).Value2 = '$/caja 50 unidades'
$ws.Cells.Item($row,15).Value2 = 'Región de Arica y Parinacota'
$ws.Cells.Item($row,16).Value2 = 230
$ws.Cells.Item($row,17).Value2 = 50
$ws.Cells.Item($row,18).Value2 = 'Hortaliza'
